$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column O, 2021 data -------------------------------------------------
# O4: header year 2021 -> same format as N4 (plain number)
$ws.Range("N4").Copy()
$ws.Range("O4").PasteSpecial(-4122)
$ws.Range("O4").Value = 2021

# O5: Total row -> same format as N5 (bold, right align) + numeric 0.0 format
$ws.Range("N5").Copy()
$ws.Range("O5").PasteSpecial(-4122)
$ws.Range("O5").NumberFormat = "0.0"
$ws.Range("O5").Value = 689

# O6-O8: same format as N10 (which already carries the 0.0 number format)
$ws.Range("N10").Copy()
$ws.Range("O6").PasteSpecial(-4122)
$ws.Range("O6").Value = 94.1

$ws.Range("O7").PasteSpecial(-4122)
$ws.Range("O7").Value = 147.1

$ws.Range("O8").PasteSpecial(-4122)
$ws.Range("O8").Value = 10.1

# O9: dash placeholder -> same base format as N9, plus 0.0 number format
$ws.Range("N9").Copy()
$ws.Range("O9").PasteSpecial(-4122)
$ws.Range("O9").NumberFormat = "0.0"
$ws.Range("O9").Value = "-"

# O10-O14: same format as N10
$ws.Range("N10").Copy()
$ws.Range("O10").PasteSpecial(-4122)
$ws.Range("O10").Value = 82.1

$ws.Range("O11").PasteSpecial(-4122)
$ws.Range("O11").Value = 145.30000000000001

$ws.Range("O12").PasteSpecial(-4122)
$ws.Range("O12").Value = 98.8

$ws.Range("O13").PasteSpecial(-4122)
$ws.Range("O13").Value = 98.7

$ws.Range("O14").PasteSpecial(-4122)
$ws.Range("O14").Value = 1.8

# O15: dash placeholder, same style as O9
$ws.Range("N15").Copy()
$ws.Range("O15").PasteSpecial(-4122)
$ws.Range("O15").NumberFormat = "0.0"
$ws.Range("O15").Value = "-"

# O16: bottom row with border, same base as N16, plus 0.0 number format
$ws.Range("N16").Copy()
$ws.Range("O16").PasteSpecial(-4122)
$ws.Range("O16").NumberFormat = "0.0"
$ws.Range("O16").Value = 10.9

$ws.Application.CutCopyMode = $false

$ws.Range("P5").Select()
